$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear all existing content/formatting in the used area and reset row heights ---
$ws.Range("A1:J40").Clear()
$ws.Rows("1:40").AutoFit()

# --- Column widths ---
$ws.Columns("F").ColumnWidth = 31.5703125
$ws.Columns("G").ColumnWidth = 89.5703125

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("E3").Select()

# --- Row 2: big note cell ---
$ws.Range("G2").Value = "Given k,
we can use equation (2) to compute L. In our case,
we chose k to be 13, and L such that the probability
of missing a neighbor within the distance of 0.2 is
less than 2.5%. The distance of 0.2 was chosen as a
reasonable estimate of the threshold when two documents are very similar. In general, this distance will
depend on the application, and Datar et al. (2004)
suggest guessing the value and then doing a binary
search to set it more accurately. We set k to 13 it achieved a reasonable balance between time
spent computing the distances and the time spent
computing the hash functions"
$ws.Range("G2").WrapText = $true
$ws.Rows(2).RowHeight = 180

# --- Rows 3-8: k / Pcoll / Pcoll power k / 1-Pcoll^k / delta / L ---
$ws.Range("E3").Value = "k"
$ws.Range("F3").Value = 16
$ws.Range("F3").Style = "Good"

$ws.Range("E4").Value = "Pcoll"
$ws.Range("F4").Value = 0.75
$ws.Range("F4").Style = "Good"

$ws.Range("E5").Value = "Pcoll power k"
$ws.Range("F5").Formula = "=POWER(F4,F3)"
$ws.Range("F5").NumberFormat = "0.0000000000"

$ws.Range("E6").Value = "1-Pcoll ^ k"
$ws.Range("F6").Formula = "=1-F5"
$ws.Range("F6").NumberFormat = "0.0000000000"

$ws.Range("E7").Value = "delta"
$ws.Range("F7").Value = 0.2
$ws.Range("F7").Style = "Good"

$ws.Range("E8").Value = "L"
$ws.Range("F8").Formula = "=LOG(F7,F6)"

# --- Row 10: max tweets ---
$ws.Range("E10").Value = "max tweets:"
$ws.Range("F10").Value = 100000
$ws.Range("G10").Value = "Rate of growth of a thread is measured by the
number of tweets that belong to that thread in a window of 100,000 tweets, starting from the beginning
of the thread."
$ws.Range("G10").WrapText = $true
$ws.Rows(10).RowHeight = 60

# --- Row 11: entropy ---
$ws.Range("E11").Value = "entropy"
$ws.Range("F11").Value = "<3.5"
$ws.Range("G11").Value = "entropy (< 3.5) to the back of the list, while we order other threads by the number of unique users.
A sign test showed this approach to be significantly better (p ≤ 0.01) than all of the previous ranking methods. Table 3 shows the effect of varying the entropy threshold at which threads are moved to the back of the list. We can see that adding information about entropy improves results regardless of the threshold we choose"
$ws.Range("G11").WrapText = $true
$ws.Rows(11).RowHeight = 90

# --- Row 12: ? / 1000-3000 / fixed number of documents ---
$ws.Range("E12").Value = "?"
$ws.Range("F12").Value = "1000-3000"
$ws.Range("G12").Value = "fixed number of most recent documents. We set this number to 2000; preliminary experiments showed that values between 1000 and 3000 all yield very similar results. "
$ws.Range("G12").WrapText = $true
$ws.Rows(12).RowHeight = 30

# --- Row 13: t / threshold for tweets closeness / In our experiments ---
$ws.Range("E13").Value = "t"
$ws.Range("F13").Value = "threshold for tweets closeness
If t is set very high, we will have few very big and broad threads, whereas setting t very low will result in many very specific and very small threads."
$ws.Range("F13").WrapText = $true
$ws.Range("G13").Value = "In our experiments, we set t = 0.5. We experimented with different values of t and found that for t ∈ [0.5,0.6] results are very much the same"
$ws.Range("G13").WrapText = $true
$ws.Rows(13).RowHeight = 90

# --- Row 14: bucket max size ---
$ws.Range("E14").Value = "bucket max size"

# --- Row 27: missing / we order the elements ---
$ws.Range("D27").Value = "missing"
$ws.Range("E27").Value = "we order the elements
of S according to the number of hash tables where
the collision occurred. We take the top 3L elements
of that ordered set and compare the new document
only to them"

# --- Row 28: missing / fastest growing? ---
$ws.Range("D28").Value = "missing"
$ws.Range("E28").Value = "fastest growing? "

# --- Rows 32-36: method description paragraphs ---
$ws.Range("F32").Value = "We first run our streaming FSD
system and assign a novelty score to each tweet"
$ws.Range("F32").WrapText = $true
$ws.Rows(32).RowHeight = 75

$ws.Range("F33").Value = "In
addition, since the score is based on a cosine distance to the nearest tweet, for each tweet we also
output which other tweet it is most similar to"
$ws.Range("F33").WrapText = $true
$ws.Rows(33).RowHeight = 105

$ws.Range("F34").Value = "links relation: tweet a links to tweet b if b is the nearest
neighbor of a and 1 −cos(a, b) < t"
$ws.Range("F34").WrapText = $true
$ws.Rows(34).RowHeight = 45

$ws.Range("F35").Value = "for each tweet a we either
assign it to an existing thread if its nearest neighbor
is within distance t, or say that a is the first tweet in
a new thread."
$ws.Range("F35").WrapText = $true
$ws.Rows(35).RowHeight = 90

$ws.Range("F36").Value = "If we assign a to an existing thread,
we assign it to the same thread to which its nearest
neighbor belongs."
$ws.Range("F36").WrapText = $true
$ws.Rows(36).RowHeight = 75

Write-Host "done"
